# Apply text edits to the TFEC sheet of the Electricity_demand workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TFEC")

# --- Text corrections (renaming of "Use" column entries) ---
$ws.Range("D5").Value  = "Residential appliances"
$ws.Range("D6").Value  = "Residential appliances new users"
$ws.Range("D7").Value  = "Commercial uses"
$ws.Range("D9").Value  = "Cars"

# --- Column D width adjustment (now best-fits the renamed/shorter labels) ---
$ws.Columns.Item(4).ColumnWidth = 27.8

# --- Update active selection ---
$ws.Range("D10").Select() | Out-Null
